# Fix typo in highlights:
#   "Determination of the heat density of decarbonized centralized heat
#    networks in 2050"
# becomes
#   "Disclosing the heat density of decarbonized centralized heat
#    networks in 2050"
#
# The paragraph's existing text is a single run immediately followed by the
# document's (hidden) "_GoBack" bookmark. We need the final OOXML to contain
# three runs - "Disclosing", " ", "the heat density..." - with the
# (collapsed) _GoBack bookmark sitting between the 2nd and 3rd runs, exactly
# like Word leaves it after an in-place edit.

$d = $word.ActiveDocument

$oldPrefix = "Determination of"
$newPrefix = "Disclosing"

# Locate the paragraph containing the text to edit.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($oldPrefix)) {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# Boundaries (within the original text) around the single space that
# separates "Determination of" from "the heat density...":
#   pStart .......... splitAfterOld   -> "Determination of"
#   splitAfterOld .... splitBeforeRest -> " "
#   splitBeforeRest .. "the heat density of decarbonized centralized heat networks in 2050"
$splitAfterOld = $pStart + $oldPrefix.Length
$splitBeforeRest = $splitAfterOld + 1

# Temporarily move the _GoBack bookmark so that it spans exactly that single
# space. Adding it in one shot (rather than as two separate collapsed
# bookmarks) makes the run split cleanly into three pieces with correct
# xml:space="preserve" handling on each fragment.
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}
$spaceRange = $d.Range($splitAfterOld, $splitBeforeRest)
$d.Bookmarks.Add("_GoBack", $spaceRange)

# Replace "Determination of" with "Disclosing" in the now-isolated first run.
$prefixRange = $d.Range($pStart, $splitAfterOld)
$prefixRange.Text = $newPrefix

# Drop the temporary bookmark and re-add _GoBack collapsed right after the
# space (i.e. between the space run and the "the heat density..." run),
# matching the structure Word produces after such an edit.
$goBack2 = $d.Bookmarks("_GoBack")
$goBack2.Delete()
$bmPos = $pStart + $newPrefix.Length + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
